# Applies the "Add files via upload" commit:
#  - profit sheet (Worksheets(1)): inserts 9 new rows (9:17) of "rate" metrics
#    with D = metric name, E = "自定义指标" (category). Existing rows 9-46
#    shift down to 18-55.
#  - balance sheet (Worksheets(2)): fills the missing E63 category cell and
#    adjusts the view/selection.
#  - CodeName properties set to match Excel's generated defaults (VBA
#    project module names).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Workbook / sheet VBA code names (cosmetic metadata, harmless to set).
# ---------------------------------------------------------------------
$wb.CodeName = "ThisWorkbook"

$codeNames = @("Sheet1","Sheet2","Sheet3","Sheet4","Sheet5","Sheet6","Sheet7","Sheet8","Sheet9","Sheet10","Sheet11")
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).CodeName = $codeNames[$i - 1]
}

# ---------------------------------------------------------------------
# "profit" sheet (index 1): insert the 9 new ratio rows at the top of the
# data block (rows 9:17), pushing the existing rows down to 18:55.
# ---------------------------------------------------------------------
$profit = $wb.Worksheets.Item(1)

$profit.Rows("9:17").Insert()

# Written in this exact order so the shared-string table picks up the new
# labels in the same sequence as the authored workbook (rows are filled
# out of left-to-right/top-to-bottom order further below).
$profit.Cells.Item(9, 4).Value  = "毛利润率[%]"
$profit.Cells.Item(10, 4).Value = "核心利润率[%]"
$profit.Cells.Item(11, 4).Value = "净利润[%]"
$profit.Cells.Item(14, 4).Value = "销售费用率[%]"
$profit.Cells.Item(16, 4).Value = "研发费用率[%]"
$profit.Cells.Item(17, 4).Value = "财务费用率[%]"
$profit.Cells.Item(12, 4).Value = "四费费率[%]"
$profit.Cells.Item(13, 4).Value = "三费费率[%]"
$profit.Cells.Item(15, 4).Value = "管理费用率[%]"

for ($r = 9; $r -le 17; $r++) {
    $profit.Cells.Item($r, 5).Value = "自定义指标"
}

# ---------------------------------------------------------------------
# "balance" sheet (index 2): row 63 was missing its category cell (E) --
# every sibling row in that block carries it.
# ---------------------------------------------------------------------
$balance = $wb.Worksheets.Item(2)
$balance.Cells.Item(63, 5).Value = "流动负债"

# ---------------------------------------------------------------------
# Selection / scroll-position bookkeeping. Select on the profit sheet
# first, then finish on the balance sheet so it remains the active tab
# (matches the workbook's activeTab).
# ---------------------------------------------------------------------
$profit.Range("G15").Select()
$balance.Range("E65").Select()
